$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextCell $ws.Range("D2") "37.495.90"
Set-TextCell $ws.Range("E2") "  +5.58%  "
Set-TextCell $ws.Range("D3") "2.054.15"
Set-TextCell $ws.Range("E3") "  +3.96%  "
Set-TextCell $ws.Range("E4") "  +0.04%  "
Set-TextCell $ws.Range("D5") "252.56"
Set-TextCell $ws.Range("E5") "  +3.37%  "
Set-TextCell $ws.Range("D6") "0.651"
Set-TextCell $ws.Range("E6") "  +2.45%  "
Set-TextCell $ws.Range("D7") "66.64"
Set-TextCell $ws.Range("E7") "  +17.30%  "
Set-TextCell $ws.Range("E8") "  +0.00%  "
Set-TextCell $ws.Range("E9") "  +6.69%  "
Set-TextCell $ws.Range("D10") "59.87"
Set-TextCell $ws.Range("E10") "  +3.70%  "
Set-TextCell $ws.Range("E11") "  +4.95%  "
Set-TextCell $ws.Range("E12") "  +1.39%  "
Set-TextCell $ws.Range("D13") "0.909"
Set-TextCell $ws.Range("E13") "  -3.56%  "
Set-TextCell $ws.Range("D14") "14.99"
Set-TextCell $ws.Range("E14") "  +4.89%  "
Set-TextCell $ws.Range("D15") "2.355.94"
Set-TextCell $ws.Range("E15") "  +4.08%  "
Set-TextCell $ws.Range("D16") "21.59"
Set-TextCell $ws.Range("E16") "  +23.08%  "
Set-TextCell $ws.Range("D17") "5.59"
Set-TextCell $ws.Range("E17") "  +6.50%  "
Set-TextCell $ws.Range("D18") "2.039.90"
Set-TextCell $ws.Range("E18") "  +3.24%  "
Set-TextCell $ws.Range("D19") "37.334.52"
Set-TextCell $ws.Range("D20") "73.70"
Set-TextCell $ws.Range("E20") "  +3.10%  "
Set-TextCell $ws.Range("D21") "0.0₃0878"
Set-TextCell $ws.Range("E21") "  +4.67%  "
Set-TextCell $ws.Range("E22") "  +6.32%  "
Set-TextCell $ws.Range("D23") "240.41"
Set-TextCell $ws.Range("E23") "  +3.56%  "
Set-TextCell $ws.Range("D24") "2.66"
Set-TextCell $ws.Range("E24") "  +3.84%  "
Set-TextCell $ws.Range("E25") "  -0.03%  "
Set-TextCell $ws.Range("D26") "2.40"
Set-TextCell $ws.Range("E26") "  +4.63%  "
Set-TextCell $ws.Range("D27") "9.78"
Set-TextCell $ws.Range("E27") "  +8.12%  "
Set-TextCell $ws.Range("D28") "160.65"
Set-TextCell $ws.Range("E28") "  -1.74%  "
Set-TextCell $ws.Range("D29") "20.07"
Set-TextCell $ws.Range("E29") "  +4.90%  "
Set-TextCell $ws.Range("B30") "Kaspa"
Set-TextCell $ws.Range("C30") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D30") "0.117"
Set-TextCell $ws.Range("E30") "  +28.04%  "
Set-TextCell $ws.Range("B31") "Filecoin"
Set-TextCell $ws.Range("C31") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D31") "5.27"
Set-TextCell $ws.Range("E31") "  +8.66%  "
Set-TextCell $ws.Range("B32") "Stellar"
Set-TextCell $ws.Range("C32") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D32") "0.122"
Set-TextCell $ws.Range("E32") "  +3.41%  "
Set-TextCell $ws.Range("E33") "  +7.95%  "
Set-TextCell $ws.Range("D34") "4.75"
Set-TextCell $ws.Range("E34") "  +11.28%  "
Set-TextCell $ws.Range("D35") "0.0623"
Set-TextCell $ws.Range("E35") "  +5.75%  "
Set-TextCell $ws.Range("E36") "  +3.70%  "
Set-TextCell $ws.Range("E37") "  +4.50%  "
Set-TextCell $ws.Range("E38") "  -0.06%  "
Set-TextCell $ws.Range("D39") "6.09"
Set-TextCell $ws.Range("E39") "  +19.45%  "
Set-TextCell $ws.Range("D40") "3.00"
Set-TextCell $ws.Range("E40") "  +34.36%  "
Set-TextCell $ws.Range("E41") "  +17.16%  "
Set-TextCell $ws.Range("E42") "  +2.74%  "
Set-TextCell $ws.Range("E43") "  +4.58%  "
Set-TextCell $ws.Range("B44") "ARBITRUM"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Range("D44") "1.15"
Set-TextCell $ws.Range("E44") "  +6.43%  "
Set-TextCell $ws.Range("B45") "VeChain"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D45") "0.0219"
Set-TextCell $ws.Range("E45") "  +4.09%  "
Set-TextCell $ws.Range("D46") "17.01"
Set-TextCell $ws.Range("E46") "  +7.44%  "
Set-TextCell $ws.Range("E47") "  +6.91%  "
Set-TextCell $ws.Range("D48") "95.88"
Set-TextCell $ws.Range("E48") "  +5.26%  "
Set-TextCell $ws.Range("D49") "1.422.32"
Set-TextCell $ws.Range("E49") "  +3.30%  "
Set-TextCell $ws.Range("E50") "  +2.24%  "
Set-TextCell $ws.Range("D51") "46.74"
Set-TextCell $ws.Range("E51") "  +1.83%  "
